$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    $val = $cell.Value2

    if ($null -eq $val) {
        continue
    }

    $text = [string]$val

    if ($text.Length -gt 0 -and $text.Substring(0,1) -eq "~") {
        $newText = $text.Substring(1)
        if ($newText.Length -gt 0 -and $newText.Substring(0,1) -eq ",") {
            $newText = $newText.Substring(1)
        }

        if ($newText -eq "") {
            $cell.ClearContents()
        } else {
            # Force text interpretation so values like " 22" are not
            # silently converted to the number 22, then restore the
            # cell's original (default) style so no stray style index
            # is left behind on the cell.
            $cell.NumberFormat = "@"
            $cell.Value2 = $newText
            $cell.Style = "Normal"
        }
    }
}
